$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows whose match data changed (F:V); A:E (Indice/pais/torneio/temporada/data_partida) are unchanged ---
# Row 3
$ws.Range("F3").Value = "Gol Gohar"
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = "Foolad"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 2.62
$ws.Range("K3").Value = "08/08/2023 06:12"
$ws.Range("L3").Value = 2.77
$ws.Range("M3").Value = "09/08/2023 17:47"
$ws.Range("N3").Value = 2.63
$ws.Range("O3").Value = "08/08/2023 06:12"
$ws.Range("P3").Value = 2.54
$ws.Range("Q3").Value = "09/08/2023 17:47"
$ws.Range("R3").Value = 2.85
$ws.Range("S3").Value = "08/08/2023 06:12"
$ws.Range("T3").Value = 3.17
$ws.Range("U3").Value = "09/08/2023 17:47"
$ws.Range("V3").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/gol-gohar-foolad/xAkf0Npf/"

# Row 4
$ws.Range("F4").Value = "Persepolis"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = "Aluminium Arak"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1.37
$ws.Range("K4").Value = "02/08/2023 06:12"
$ws.Range("L4").Value = 1.46
$ws.Range("M4").Value = "09/08/2023 16:33"
$ws.Range("N4").Value = 3.83
$ws.Range("O4").Value = "02/08/2023 06:12"
$ws.Range("P4").Value = 3.63
$ws.Range("Q4").Value = "09/08/2023 16:33"
$ws.Range("R4").Value = 8.02
$ws.Range("S4").Value = "02/08/2023 06:12"
$ws.Range("T4").Value = 9.26
$ws.Range("U4").Value = "09/08/2023 16:33"
$ws.Range("V4").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/persepolis-aluminium-arak/xpKbcLUC/"

# Row 5
$ws.Range("F5").Value = "Malavan"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = "Zob Ahan"
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 2.62
$ws.Range("K5").Value = "08/08/2023 06:12"
$ws.Range("L5").Value = 3.15
$ws.Range("M5").Value = "09/08/2023 16:52"
$ws.Range("N5").Value = 2.62
$ws.Range("O5").Value = "08/08/2023 06:12"
$ws.Range("P5").Value = 2.43
$ws.Range("Q5").Value = "09/08/2023 16:52"
$ws.Range("R5").Value = 2.86
$ws.Range("S5").Value = "08/08/2023 06:12"
$ws.Range("T5").Value = 2.92
$ws.Range("U5").Value = "09/08/2023 16:52"
$ws.Range("V5").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/malavan-zob-ahan/4din2qGs/"

# Row 35
$ws.Range("F35").Value = "Aluminium Arak"
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = "Tractor"
$ws.Range("I35").Value = 4
$ws.Range("J35").Value = 2.83
$ws.Range("K35").Value = "05/10/2023 09:42"
$ws.Range("L35").Value = 3.18
$ws.Range("M35").Value = "05/10/2023 15:43"
$ws.Range("N35").Value = 2.74
$ws.Range("O35").Value = "05/10/2023 09:42"
$ws.Range("P35").Value = 2.45
$ws.Range("Q35").Value = "05/10/2023 15:43"
$ws.Range("R35").Value = 2.71
$ws.Range("S35").Value = "05/10/2023 09:42"
$ws.Range("T35").Value = 2.87
$ws.Range("U35").Value = "05/10/2023 15:43"
$ws.Range("V35").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/aluminium-arak-tractor/466ZSTPT/"

# Row 36
$ws.Range("F36").Value = "Shams Azar Qazvin"
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = "Zob Ahan"
$ws.Range("I36").Value = 2
$ws.Range("J36").Value = 2.98
$ws.Range("K36").Value = "05/10/2023 09:42"
$ws.Range("L36").Value = 2.83
$ws.Range("M36").Value = "05/10/2023 15:32"
$ws.Range("N36").Value = 2.78
$ws.Range("O36").Value = "05/10/2023 09:42"
$ws.Range("P36").Value = 2.65
$ws.Range("Q36").Value = "05/10/2023 15:16"
$ws.Range("R36").Value = 2.61
$ws.Range("S36").Value = "05/10/2023 09:42"
$ws.Range("T36").Value = 2.93
$ws.Range("U36").Value = "05/10/2023 15:32"
$ws.Range("V36").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/shams-azar-qazvin-zob-ahan/Aq7VT9AN/"

# Row 50
$ws.Range("F50").Value = "Esteghlal F.C."
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = "Aluminium Arak"
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 1.49
$ws.Range("K50").Value = "26/10/2023 07:42"
$ws.Range("L50").Value = 1.54
$ws.Range("M50").Value = "27/10/2023 16:24"
$ws.Range("N50").Value = 3.49
$ws.Range("O50").Value = "26/10/2023 07:42"
$ws.Range("P50").Value = 3.42
$ws.Range("Q50").Value = "27/10/2023 16:25"
$ws.Range("R50").Value = 6.44
$ws.Range("S50").Value = "26/10/2023 07:42"
$ws.Range("T50").Value = 7.84
$ws.Range("U50").Value = "27/10/2023 16:25"
$ws.Range("V50").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/esteghlal-teh-aluminium-arak/lzS4r9PG/"

# Row 51
$ws.Range("F51").Value = "Zob Ahan"
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = "Tractor"
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3.1
$ws.Range("K51").Value = "26/10/2023 07:42"
$ws.Range("L51").Value = 3.13
$ws.Range("M51").Value = "27/10/2023 16:17"
$ws.Range("N51").Value = 2.74
$ws.Range("O51").Value = "26/10/2023 07:42"
$ws.Range("P51").Value = 2.69
$ws.Range("Q51").Value = "27/10/2023 16:18"
$ws.Range("R51").Value = 2.35
$ws.Range("S51").Value = "26/10/2023 07:42"
$ws.Range("T51").Value = 2.63
$ws.Range("U51").Value = "27/10/2023 16:18"
$ws.Range("V51").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/zob-ahan-tractor/I9W8sTvN/"

# Row 58
$ws.Range("F58").Value = "Tractor"
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = "Nassaji Mazandaran"
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1.58
$ws.Range("K58").Value = "01/11/2023 01:12"
$ws.Range("L58").Value = 1.71
$ws.Range("M58").Value = "02/11/2023 12:58"
$ws.Range("N58").Value = 3.35
$ws.Range("O58").Value = "01/11/2023 01:12"
$ws.Range("P58").Value = 3.26
$ws.Range("Q58").Value = "02/11/2023 12:58"
$ws.Range("R58").Value = 5.39
$ws.Range("S58").Value = "01/11/2023 01:12"
$ws.Range("T58").Value = 5.53
$ws.Range("U58").Value = "02/11/2023 12:58"
$ws.Range("V58").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/tractor-mazandaran/lEbFlpnj/"

# Row 59
$ws.Range("F59").Value = "Shams Azar Qazvin"
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = "Esteghlal F.C."
$ws.Range("I59").Value = 2
$ws.Range("J59").Value = 5.02
$ws.Range("K59").Value = "01/11/2023 00:42"
$ws.Range("L59").Value = 5.08
$ws.Range("M59").Value = "02/11/2023 12:59"
$ws.Range("N59").Value = 3.4
$ws.Range("O59").Value = "01/11/2023 00:42"
$ws.Range("P59").Value = 2.99
$ws.Range("Q59").Value = "02/11/2023 12:59"
$ws.Range("R59").Value = 1.63
$ws.Range("S59").Value = "01/11/2023 00:42"
$ws.Range("T59").Value = 1.85
$ws.Range("U59").Value = "02/11/2023 12:59"
$ws.Range("V59").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/shams-azar-qazvin-esteghlal-teh/StFWd8vA/"

# Row 62
$ws.Range("F62").Value = "Havadar SC"
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = "Paykan"
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2.17
$ws.Range("K62").Value = "02/11/2023 00:42"
$ws.Range("L62").Value = 2.15
$ws.Range("M62").Value = "03/11/2023 12:28"
$ws.Range("N62").Value = 2.65
$ws.Range("O62").Value = "02/11/2023 00:42"
$ws.Range("P62").Value = 2.41
$ws.Range("Q62").Value = "03/11/2023 12:28"
$ws.Range("R62").Value = 3.62
$ws.Range("S62").Value = "02/11/2023 00:42"
$ws.Range("T62").Value = 4.31
$ws.Range("U62").Value = "03/11/2023 12:28"
$ws.Range("V62").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/havadar-sc-paykan/KvfBkQWq/"

# Row 63
$ws.Range("F63").Value = "Mes Rafsanjan"
$ws.Range("G63").Value = 3
$ws.Range("H63").Value = "Foolad"
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 2.2
$ws.Range("K63").Value = "02/11/2023 00:42"
$ws.Range("L63").Value = 2.3
$ws.Range("M63").Value = "03/11/2023 12:29"
$ws.Range("N63").Value = 2.65
$ws.Range("O63").Value = "02/11/2023 00:42"
$ws.Range("P63").Value = 2.48
$ws.Range("Q63").Value = "03/11/2023 12:29"
$ws.Range("R63").Value = 3.55
$ws.Range("S63").Value = "02/11/2023 00:42"
$ws.Range("T63").Value = 4.29
$ws.Range("U63").Value = "03/11/2023 12:29"
$ws.Range("V63").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/mes-rafsanjan-foolad/z7Iven9M/"

# --- Append new rows 74-76 ---
# Copy the A/E column formatting (bold+border style, date number format) from row 73, the last existing row
$ws.Range("A73").Copy()
$ws.Range("A74").PasteSpecial(-4122)
$ws.Range("E73").Copy()
$ws.Range("E74").PasteSpecial(-4122)
$ws.Range("A73").Copy()
$ws.Range("A75").PasteSpecial(-4122)
$ws.Range("E73").Copy()
$ws.Range("E75").PasteSpecial(-4122)
$ws.Range("A73").Copy()
$ws.Range("A76").PasteSpecial(-4122)
$ws.Range("E73").Copy()
$ws.Range("E76").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 74
$ws.Range("A74").Value = 73
$ws.Range("B74").Value = "iran"
$ws.Range("C74").Value = "persian-gulf-pro-league"
$ws.Range("D74").Value = "2023-2024"
$ws.Range("E74").Value = 45254.52083333334
$ws.Range("F74").Value = "Aluminium Arak"
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = "Paykan"
$ws.Range("I74").Value = 1
$ws.Range("J74").Value = 1.85
$ws.Range("K74").Value = "22/11/2023 15:12"
$ws.Range("L74").Value = 2.1
$ws.Range("M74").Value = "24/11/2023 12:10"
$ws.Range("N74").Value = 2.8
$ws.Range("O74").Value = "22/11/2023 15:12"
$ws.Range("P74").Value = 2.36
$ws.Range("Q74").Value = "24/11/2023 12:10"
$ws.Range("R74").Value = 4.58
$ws.Range("S74").Value = "22/11/2023 15:12"
$ws.Range("T74").Value = 4.43
$ws.Range("U74").Value = "24/11/2023 12:10"
$ws.Range("V74").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/aluminium-arak-paykan/2o67LExL/"

# Row 75
$ws.Range("A75").Value = 74
$ws.Range("B75").Value = "iran"
$ws.Range("C75").Value = "persian-gulf-pro-league"
$ws.Range("D75").Value = "2023-2024"
$ws.Range("E75").Value = 45254.52083333334
$ws.Range("F75").Value = "Havadar SC"
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = "Malavan"
$ws.Range("I75").Value = 1
$ws.Range("J75").Value = 2.51
$ws.Range("K75").Value = "22/11/2023 15:12"
$ws.Range("L75").Value = 3.16
$ws.Range("M75").Value = "24/11/2023 11:20"
$ws.Range("N75").Value = 2.65
$ws.Range("O75").Value = "22/11/2023 15:12"
$ws.Range("P75").Value = 2.5
$ws.Range("Q75").Value = "24/11/2023 11:20"
$ws.Range("R75").Value = 2.96
$ws.Range("S75").Value = "22/11/2023 15:12"
$ws.Range("T75").Value = 2.81
$ws.Range("U75").Value = "24/11/2023 11:33"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/havadar-sc-malavan/tx5BKYiR/"

# Row 76
$ws.Range("A76").Value = 75
$ws.Range("B76").Value = "iran"
$ws.Range("C76").Value = "persian-gulf-pro-league"
$ws.Range("D76").Value = "2023-2024"
$ws.Range("E76").Value = 45254.52083333334
$ws.Range("F76").Value = "Mes Rafsanjan"
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = "Sanat Naft"
$ws.Range("I76").Value = 2
$ws.Range("J76").Value = 1.53
$ws.Range("K76").Value = "22/11/2023 15:12"
$ws.Range("L76").Value = 1.58
$ws.Range("M76").Value = "24/11/2023 12:25"
$ws.Range("N76").Value = 3.43
$ws.Range("O76").Value = "22/11/2023 15:12"
$ws.Range("P76").Value = 3.36
$ws.Range("Q76").Value = "24/11/2023 12:25"
$ws.Range("R76").Value = 6.31
$ws.Range("S76").Value = "22/11/2023 15:12"
$ws.Range("T76").Value = 7.1
$ws.Range("U76").Value = "24/11/2023 12:25"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/mes-rafsanjan-sanat-naft/MVnvGh0r/"
